$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.243.88"
$ws.Range("E2").Value = "  -3.35%  "
$ws.Range("D3").Value = "3.158.64"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.98%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.156.60"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("E9").Value = "  -3.10%  "
$ws.Range("E10").Value = "  -6.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("E12").Value = "  -4.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.64%  "
$ws.Range("D15").Value = "3.677.48"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").Value = "64.301.66"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "3.161.59"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("E19").Value = "  -4.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.713"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.27%  "
$ws.Range("E30").Value = "  -27.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.83%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.67%  "
$ws.Range("E35").Value = "  -5.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "0.0₃0727"
$ws.Range("E38").Value = "  -7.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "451.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.13%  "
$ws.Range("E40").Value = "  -5.27%  "
$ws.Range("E41").Value = "  -4.72%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.119"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.23%  "
$ws.Range("D44").Value = "2.853.11"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.269"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.95%  "
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.22%  "
